$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-18 from
# 2023-10-13 (45212) to 2023-10-22 (45221)
$newDate = Get-Date -Year 2023 -Month 10 -Day 22 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

for ($row = 2; $row -le 18; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
